{"js": "// Office.js (Word JavaScript API) implementation of the journal update.\n//\n// The underlying OOXML diff shows five edits in the document:\n//   1) The title-page author line \"Polina Prinii\" gets spell-check\n//      proofing marks (w:proofErr) added around \"Prinii\" \u2014 a pure\n//      run-split with no visible text change.\n//   2) The \"virtualenv\" mention gets the same proofing-mark treatment\n//      \u2014 again no visible text change.\n//   3) \"Work commenced at 2:40pm.\" is re-merged from three runs back\n//      into a single run \u2014 no visible text change.\n//   4) \"Jupyter\" gets proofing marks added \u2014 no visible text change.\n//   5) New journal content is appended after the \"Additionally, found\n//      that PyCharm ...\" paragraph: two more sentences are appended to\n//      that paragraph, followed by a blank spacer paragraph, a new\n//      paragraph (\"For this reason I will just perform ...\"), and two\n//      more blank spacer paragraphs.\n//\n// (1)-(4) are Word's own automatic spell-check bookkeeping (added by the\n// desktop editor as a side effect of someone typing/correcting text) and\n// do not change the document's visible text \u2014 they are not something an\n// automation script would normally (or even typically can) reproduce.\n// The only substantive, author-driven content edit is (5), so that is\n// what this script performs.\n\nconst body = context.document.body;\n\n// Locate the paragraph that ends the \"Additionally, found ...\" sentence\n// by searching for its distinctive text, rather than relying on a\n// paragraph index (more robust to unrelated layout differences).\nconst searchText =\n  \"Additionally, found that PyCharm for some odd reason limits the print \" +\n  \"output of a data-frame to two columns and 10 rows. To get around this \" +\n  \"I must use pd options around display to show more than 10 rows and 2 \" +\n  \"columns at a time.\";\n\nconst results = body.search(searchText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the \"Additionally, found ...\" paragraph.');\n}\n\nconst found = results.items[0];\nconst targetParagraphs = found.paragraphs;\ntargetParagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = targetParagraphs.items[0];\n\n// 1) Append the two new sentences to the existing paragraph.\ntargetParagraph.insertText(\n  \" Having applied the pd options to expand the display of columns and rows\" +\n    \", it would expand the view of the columns and not the rows. Unsure if \" +\n    \"there is some issue around the rows but the pd option won\\u2019t apply to rows.\",\n  \"End\"\n);\nawait context.sync();\n\n// 2) Insert a blank \"NoSpacing\" spacer paragraph right after it. Leaving\n//    styleBuiltIn untouched lets the new paragraph inherit the same\n//    paragraph mark formatting (pStyle \"NoSpacing\" + Emphasis rPr) as its\n//    neighbour, matching how Word itself extends a run of like-styled\n//    paragraphs.\nconst spacer1 = targetParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\n\n// 3) Insert the new \"For this reason I ...\" paragraph after the spacer.\nconst reasonParagraph = spacer1.insertParagraph(\n  \"For this reason I will just perform a number of print functions to \" +\n    \"validate if the data came through okay, such as print columns \" +\n    \"values, to show columns names, a count function to see how many \" +\n    \"rows and so on.\",\n  \"After\"\n);\nawait context.sync();\n\n// 4) Insert two more blank spacer paragraphs after that.\nconst spacer2 = reasonParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\n\nconst spacer3 = spacer2.insertParagraph(\"\", \"After\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) implementation of the journal update.\n#\n# The underlying OOXML diff shows five edits in the document:\n#   1) The title-page author line \"Polina Prinii\" gets spell-check\n#      proofing marks (w:proofErr) added around \"Prinii\" - a pure\n#      run-split with no visible text change.\n#   2) The \"virtualenv\" mention gets the same proofing-mark treatment\n#      - again no visible text change.\n#   3) \"Work commenced at 2:40pm.\" is re-merged from three runs back\n#      into a single run - no visible text change.\n#   4) \"Jupyter\" gets proofing marks added - no visible text change.\n#   5) New journal content is appended after the \"Additionally, found\n#      that PyCharm ...\" paragraph: two more sentences are appended to\n#      that paragraph, followed by a blank spacer paragraph, a new\n#      paragraph (\"For this reason I will just perform ...\"), and two\n#      more blank spacer paragraphs.\n#\n# (1)-(4) are Word's own automatic spell-check bookkeeping (added by the\n# desktop editor as a side effect of someone typing/correcting text) and\n# do not change the document's visible text - they are not something an\n# automation script would normally (or even typically can) reproduce.\n# The only substantive, author-driven content edit is (5), so that is\n# what this script performs.\n\n$d = $word.ActiveDocument\n$apos = [char]8217   # U+2019 RIGHT SINGLE QUOTATION MARK, matches the doc's curly apostrophe\n\n# Locate the paragraph that ends the \"Additionally, found ...\" sentence by\n# searching for its distinctive text, rather than relying on a paragraph\n# index (more robust to unrelated layout differences).\n$searchText = \"Additionally, found that PyCharm for some odd reason limits the print output of a data-frame to two columns and 10 rows. To get around this I must use pd options around display to show more than 10 rows and 2 columns at a time.\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\n$rng.Find.MatchCase = $true\n$found = $rng.Find.Execute($searchText)\n\nif (-not $found) {\n    throw \"Could not find the 'Additionally, found ...' paragraph.\"\n}\n\n# $rng now spans exactly the found text; get its containing paragraph.\n$targetParagraph = $rng.Paragraphs.Item(1)\n\n# 1) Append the two new sentences to the existing paragraph.\n$targetParagraph.Range.InsertAfter(\" Having applied the pd options to expand the display of columns and rows, it would expand the view of the columns and not the rows. Unsure if there is some issue around the rows but the pd option won\" + $apos + \"t apply to rows.\")\n\n$targetIndex = $targetParagraph.Range.Paragraphs.Item(1).Index\n\n# 2) Insert a blank \"NoSpacing\" spacer paragraph right after it. (Leaving\n#    its formatting alone lets it inherit the same paragraph mark\n#    formatting as its neighbour, matching Word's own behaviour when\n#    extending a run of like-styled paragraphs.)\n$d.Paragraphs.Item($targetIndex).Range.InsertParagraphAfter() | Out-Null\n\n# 3) Insert the new \"For this reason I ...\" paragraph after the spacer.\n$d.Paragraphs.Item($targetIndex + 1).Range.InsertParagraphAfter() | Out-Null\n$d.Paragraphs.Item($targetIndex + 2).Range.Text = \"For this reason I will just perform a number of print functions to validate if the data came through okay, such as print columns values, to show columns names, a count function to see how many rows and so on.\"\n\n# 4) Insert two more blank spacer paragraphs after that.\n$d.Paragraphs.Item($targetIndex + 2).Range.InsertParagraphAfter() | Out-Null\n$d.Paragraphs.Item($targetIndex + 3).Range.InsertParagraphAfter() | Out-Null\n"}
